# ------------------------------------------------------------------
# Applies the "Add files via upload" edit to Sheet1 of the workbook:
#   1. Adds a new column V mirroring column B (station header + codes)
#      for rows 26-45, with a merged V26:V27 header cell.
#   2. Fills in previously-empty Trend/sigma_trend (T/U) columns for
#      rows 28-45 with the measured values (yellow fill), using the
#      red "-" placeholder fill for the two stations without data
#      (rows 32 and 38).
#   3. Updates the active selection to X32 (best-effort scroll to M19).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. New column V = copy of column B (header + station codes) ---
$ws.Range("B26:B45").Copy($ws.Range("V26:V45"))
$ws.Range("V26:V27").Merge()

# --- 2. Fill in the Trend (T) / sigma_trend (U) columns ---------------
# Stations with measured trend values -> yellow highlight
$ws.Range("T28").Value = 0.73
$ws.Range("U28").Value = 1.6
$ws.Range("T28:U28").Interior.Color = 65535

$ws.Range("T29").Value = 2.86
$ws.Range("U29").Value = 0.8
$ws.Range("T29:U29").Interior.Color = 65535

$ws.Range("T30").Value = 2.25
$ws.Range("U30").Value = 0.8
$ws.Range("T30:U30").Interior.Color = 65535

$ws.Range("T31").Value = 37.94
$ws.Range("U31").Value = 4.1
$ws.Range("T31:U31").Interior.Color = 65535

# Station without data -> red "-" placeholder
$ws.Range("T32").Value = "-"
$ws.Range("U32").Value = "-"
$ws.Range("T32:U32").Interior.Color = 255

$ws.Range("T33").Value = -10.85
$ws.Range("U33").Value = 1.6
$ws.Range("T33:U33").Interior.Color = 65535

$ws.Range("T34").Value = 3.63
$ws.Range("U34").Value = 0.9
$ws.Range("T34:U34").Interior.Color = 65535

$ws.Range("T35").Value = 3.73
$ws.Range("U35").Value = 0.8
$ws.Range("T35:U35").Interior.Color = 65535

$ws.Range("T36").Value = 0.72
$ws.Range("U36").Value = 1.9
$ws.Range("T36:U36").Interior.Color = 65535

$ws.Range("T37").Value = -0.77
$ws.Range("U37").Value = 1.1
$ws.Range("T37:U37").Interior.Color = 65535

# Station without data -> red "-" placeholder
$ws.Range("T38").Value = "-"
$ws.Range("U38").Value = "-"
$ws.Range("T38:U38").Interior.Color = 255

$ws.Range("T39").Value = -1.72
$ws.Range("U39").Value = 1.9
$ws.Range("T39:U39").Interior.Color = 65535

$ws.Range("T40").Value = 3.84
$ws.Range("U40").Value = 0.6
$ws.Range("T40:U40").Interior.Color = 65535

$ws.Range("T41").Value = 4.56
$ws.Range("U41").Value = 2.5
$ws.Range("T41:U41").Interior.Color = 65535

$ws.Range("T42").Value = 2.38
$ws.Range("U42").Value = 2.4
$ws.Range("T42:U42").Interior.Color = 65535

$ws.Range("T43").Value = 0.16
$ws.Range("U43").Value = 1
$ws.Range("T43:U43").Interior.Color = 65535

$ws.Range("T44").Value = -6.36
$ws.Range("U44").Value = 2.6
$ws.Range("T44:U44").Interior.Color = 65535

$ws.Range("T45").Value = -6.47
$ws.Range("U45").Value = 2.1
$ws.Range("T45:U45").Interior.Color = 65535

# --- 3. View state: scroll to M19, select X32 (best effort) -----------
$ws.Range("M19").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 13
$ws.Range("X32").Select()
